# Calibrate passive detection trend
# Adds a new "parameter" row (passive_detection_past_frac) to the
# "constant" worksheet, right after the existing passive_detection_shape
# row (row 51), mirroring the layout of the neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

$newRow = 52

$ws.Cells.Item($newRow, 1).Value = "passive_detection_past_frac"
$ws.Cells.Item($newRow, 2).Value = 0.75
$ws.Cells.Item($newRow, 3).Value = "uniform"
$ws.Cells.Item($newRow, 4).Value = 0.5
$ws.Cells.Item($newRow, 5).Value = 1
$ws.Cells.Item($newRow, 7).Value = "Past passive detection rate, as a fraction of the current one"

# Match the author's updated viewport / selection state.
$ws.Activate()
$ws.Range("F50").Select()
